$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.573.89'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '1.665.72'
$ws.Range("E3").Value = '  +1.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9984'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.62'
$ws.Range("E5").Value = '  +1.53%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4621'
$ws.Range("E7").Value = '  -2.99%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2581'
$ws.Range("E8").Value = '  -0.37%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06142'
$ws.Range("E9").Value = '  +0.63%  '
$ws.Range("D10").Value = '1.663.26'
$ws.Range("E10").Value = '  +1.05%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06927'
$ws.Range("E11").Value = '  -1.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.83'
$ws.Range("E12").Value = '  +2.05%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.344'
$ws.Range("E13").Value = '  +0.00%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '75.20'
$ws.Range("E14").Value = '  +2.18%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5727'
$ws.Range("E15").Value = '  -2.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9991'
$ws.Range("E16").Value = '  -0.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.9996'
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("D18").Value = '25.584.68'
$ws.Range("E18").Value = '  +2.42%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000006693'
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D21").Value = '1.876.84'
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.431'
$ws.Range("E22").Value = '  +3.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.609'
$ws.Range("E23").Value = '  +0.55%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.220'
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '134.15'
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '14.93'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.375'
$ws.Range("E27").Value = '  -0.53%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.721'
$ws.Range("E28").Value = '  +5.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '104.25'
$ws.Range("E29").Value = '  +0.84%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.942'
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.07662'
$ws.Range("E31").Value = '  +0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.595'
$ws.Range("E32").Value = '  +0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.04340'
$ws.Range("E33").Value = '  +1.39%  '
$ws.Range("E34").Value = '  +1.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6069'
$ws.Range("E35").Value = '  +2.98%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9391'
$ws.Range("E36").Value = '  +1.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9239'
$ws.Range("E37").Value = '  +5.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.436'
$ws.Range("E38").Value = '  -5.59%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '107.31'
$ws.Range("E39").Value = '  +8.65%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9994'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.832'
$ws.Range("E41").Value = '  +4.30%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01450'
$ws.Range("E42").Value = '  -3.58%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.059'
$ws.Range("E43").Value = '  +8.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.3711'
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1109'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("E46").Value = '  +1.16%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.105'
$ws.Range("E47").Value = '  +0.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '30.98'
$ws.Range("E48").Value = '  +7.50%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.620'
$ws.Range("E49").Value = '  +7.25%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.001'
$ws.Range("E50").Value = '  +0.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.206'
$ws.Range("E51").Value = '  +3.15%  '
